# Updated capital structure database
# Applies per-cell value updates to sheet1 (France - Financial Svcs. Non-bank & Insurance)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2 ---
$ws.Range("D2").Value = 0.00289
$ws.Range("E2").Value = -0.09900000000000002
$ws.Range("F2").Value = -0.0649
$ws.Range("G2").Value = 0.002213300012862682
$ws.Range("H2").Value = 0.002213300012862682
$ws.Range("I2").Value = 0.004180278552757799
$ws.Range("J2").Value = 0.003581993926172231
$ws.Range("K2").Value = 401.54
$ws.Range("L2").Value = 0.03689200859962147
$ws.Range("M2").Value = 24.52
$ws.Range("N2").Value = 0.001800135083546237
$ws.Range("O2").Value = 0.06106490013448225
$ws.Range("P2").Value = 24.52
$ws.Range("Q2").Value = 0.001800135083546237
$ws.Range("R2").Value = 0.06106490013448225
$ws.Range("S2").Value = 0
$ws.Range("T2").Value = 0
$ws.Range("U2").Value = 36349.92
$ws.Range("V2").Value = 2.668628314685931
$ws.Range("W2").Value = 0.04411819216744908
$ws.Range("X2").Value = 0.02227973671115513
$ws.Range("Y2").Value = 0.02183845545629395
$ws.Range("Z2").Value = 0.2236274125595279
$ws.Range("AA2").Value = 0.1422133769279977
$ws.Range("AB2").Value = 0.02090303043734218
$ws.Range("AC2").Value = 0.123516859075864
$ws.Range("AD2").Value = 46077.43
$ws.Range("AE2").Value = 1.35506088036788
$ws.Range("AF2").Value = 46078.78506088037
$ws.Range("AG2").Value = 9728.865060880373
$ws.Range("AH2").Value = 0.771839138885722
$ws.Range("AI2").Value = 0.6441030139331785
$ws.Range("AJ2").Value = 0.4166525889976927
$ws.Range("AK2").Value = 0.2764703278786002
$ws.Range("AL2").Value = 2.336
$ws.Range("AM2").Value = 2.226
$ws.Range("AN2").Value = 925.2495983935744
$ws.Range("AO2").Value = 19.04965753424658
$ws.Range("AP2").Value = 195.3587361622565
$ws.Range("AQ2").Value = 19.99101527403414

# --- Row 3 ---
$ws.Range("B3").Value = "Advenis SA (ENXTPA:ADV)"
$ws.Range("D3").Value = 0.00289
$ws.Range("G3").Value = 0.03426423200859291
$ws.Range("H3").Value = 0.03426423200859291
$ws.Range("I3").Value = 0.04726100966702471
$ws.Range("J3").Value = 0.04726100966702471
$ws.Range("K3").Value = -2.96
$ws.Range("L3").Value = -0.0317937701396348
$ws.Range("M3").Value = -0
$ws.Range("N3").Value = -0
$ws.Range("O3").Value = 0
$ws.Range("P3").Value = -0
$ws.Range("Q3").Value = -0
$ws.Range("R3").Value = 0
$ws.Range("S3").Value = 0
$ws.Range("U3").Value = 7.32
$ws.Range("V3").Value = 0.2407894736842106
$ws.Range("W3").Value = -0.4892561983471074
$ws.Range("X3").Value = 0.02451430917896858
$ws.Range("Y3").Value = -0.513770507526076
$ws.Range("Z3").Value = 14.56964006259781
$ws.Range("AA3").Value = 0.6885758998435054
$ws.Range("AB3").Value = 0.02463319843993352
$ws.Range("AC3").Value = 0.6639427014035719
$ws.Range("AD3").Value = 33
$ws.Range("AF3").Value = 33
$ws.Range("AG3").Value = 25.68
$ws.Range("AH3").Value = 0.5205047318611987
$ws.Range("AI3").Value = 0.9171762090050029
$ws.Range("AJ3").Value = 0.4579172610556348
$ws.Range("AK3").Value = 0.8960223307745987
$ws.Range("AL3").Value = 2.3
$ws.Range("AM3").Value = 2.19
$ws.Range("AN3").Value = 4.269081500646831
$ws.Range("AO3").Value = 1.91304347826087
$ws.Range("AP3").Value = 3.322121604139715
$ws.Range("AQ3").Value = 2.009132420091325
$ws.Range("E3").ClearContents()
$ws.Range("T3").ClearContents()

# --- Row 4 ---
$ws.Range("B4").Value = "ABC arbitrage SA (ENXTPA:ABCA)"
$ws.Range("D4").Value = 0.0882
$ws.Range("E4").Value = 0.149
$ws.Range("G4").Value = 0.2886740331491712
$ws.Range("H4").Value = 0.2886740331491712
$ws.Range("I4").Value = 0.5538674033149171
$ws.Range("J4").Value = 0.552707182320442
$ws.Range("K4").Value = 40
$ws.Range("L4").Value = 0.5524861878453038
$ws.Range("M4").Value = 21.6
$ws.Range("N4").Value = 0.04179566563467493
$ws.Range("O4").Value = 0.54
$ws.Range("P4").Value = 21.6
$ws.Range("Q4").Value = 0.04179566563467493
$ws.Range("R4").Value = 0.54
$ws.Range("S4").Value = 0
$ws.Range("T4").Value = 0
$ws.Range("U4").Value = 14.6
$ws.Range("V4").Value = 0.02825077399380805
$ws.Range("W4").Value = 0.2489110143123833
$ws.Range("X4").Value = 0.01786275174616176
$ws.Range("Y4").Value = 0.2310482625662216
$ws.Range("Z4").Value = 0.5146065818466131
$ws.Range("AA4").Value = 0.2844267538559955
$ws.Range("AB4").Value = 0.01785947152720309
$ws.Range("AC4").Value = 0.2665672823287924
$ws.Range("AD4").Value = 1.93
$ws.Range("AE4").Value = 0
$ws.Range("AF4").Value = 1.93
$ws.Range("AG4").Value = -12.67
$ws.Range("AH4").Value = 0.003720625373508377
$ws.Range("AI4").Value = 0.01083478358502218
$ws.Range("AJ4").Value = -0.0251324063237657
$ws.Range("AK4").Value = -0.07747813856784688
$ws.Range("AL4").Value = 0.036
$ws.Range("AM4").Value = 0.036
$ws.Range("AN4").Value = 0.04730392156862745
$ws.Range("AO4").Value = 1113.888888888889
$ws.Range("AP4").Value = -0.3105392156862745
$ws.Range("AQ4").Value = 1113.888888888889

# --- Row 5 ---
$ws.Range("B5").Value = "Rothschild & Co SCA (ENXTPA:ROTH)"
$ws.Range("I5").Value = 0.0004723570021875379
$ws.Range("J5").Value = 0.0003982898394493709
$ws.Range("K5").Value = 189.5
$ws.Range("L5").Value = 0.08960234526455151
$ws.Range("M5").Value = 2.92
$ws.Range("N5").Value = 0.001254511084378759
$ws.Range("O5").Value = 0.0154089709762533
$ws.Range("P5").Value = 2.92
$ws.Range("Q5").Value = 0.001254511084378759
$ws.Range("R5").Value = 0.0154089709762533
$ws.Range("T5").Value = 0
$ws.Range("U5").Value = 4388.1
$ws.Range("V5").Value = 1.885246605946039
$ws.Range("W5").Value = 0.07992071190586648
$ws.Range("X5").Value = 0.02004516424334168
$ws.Range("Y5").Value = 0.05987554766252481
$ws.Range("Z5").Value = -0.7463795512123019
$ws.Range("AA5").Value = -0.0002972753916206412
$ws.Range("AB5").Value = 0.01923628878544371
$ws.Range("AC5").Value = -0.01953356417706436
$ws.Range("AD5").Value = 833.5
$ws.Range("AE5").Value = 1.35506088036788
$ws.Range("AF5").Value = 834.8550608803679
$ws.Range("AG5").Value = -3553.244939119632
$ws.Range("AH5").Value = 0.263989541292631
$ws.Range("AI5").Value = 0.2263734365247596
$ws.Range("AJ5").Value = 2.899081802330037
$ws.Range("AK5").Value = 5.075013387353066
$ws.Range("AN5").Value = 656.2992125984252
$ws.Range("AP5").Value = -2797.83066072412
$ws.Range("D5").ClearContents()
$ws.Range("E5").ClearContents()
$ws.Range("F5").ClearContents()

# --- Row 6 ---
$ws.Range("B6").Value = "Natixis S.A. (ENXTPA:KN)"
$ws.Range("D6").Value = -0.019
$ws.Range("E6").Value = -0.347
$ws.Range("F6").Value = -0.0649
$ws.Range("G6").Value = 0
$ws.Range("H6").Value = 0
$ws.Range("I6").Value = 0
$ws.Range("J6").Value = 0
$ws.Range("K6").Value = 175
$ws.Range("L6").Value = 0.02033984983379437
$ws.Range("O6").Value = -0
$ws.Range("R6").Value = -0
$ws.Range("U6").Value = 31939.9
$ws.Range("V6").Value = 2.972148812625624
$ws.Range("W6").Value = 0.008315672429031676
$ws.Range("X6").Value = 0.04370652134514222
$ws.Range("Y6").Value = -0.03539084891611054
$ws.Range("Z6").Value = 0.1675272987834322
$ws.Range("AA6").Value = 0
$ws.Range("AB6").Value = 0.02256977208924065
$ws.Range("AC6").Value = -0.02256977208924065
$ws.Range("AD6").Value = 45209
$ws.Range("AF6").Value = 45209
$ws.Range("AG6").Value = 13269.1
$ws.Range("AH6").Value = 0.8079470435382465
$ws.Range("AI6").Value = 0.6684023927590358
$ws.Range("AJ6").Value = 0.5525223293289749
$ws.Range("AK6").Value = 0.3717095034666293
$ws.Range("AL6").Value = 0
$ws.Range("AM6").Value = 0
$ws.Range("AN6").ClearContents()
$ws.Range("AO6").ClearContents()
$ws.Range("AP6").ClearContents()
$ws.Range("AQ6").ClearContents()

